# Apply edit described by commit: dynamic start_urls parsing, robust max_page
# handling, and improved logging for the ZorgkaartScrapy spider. As a result the
# scraped dataset now contains additional "huisartsenpraktijk" entries interleaved
# with the existing "tandartsenpraktijk" entries (rows 22-121, with the sheet
# dimension growing from A1:C61 to A1:C121).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 100,3
$data[0,0] = "huisartsenpraktijk"
$data[0,1] = "Snipmeister"
$data[0,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-snipmeister-zeist-10021848"
$data[1,0] = "huisartsenpraktijk"
$data[1,1] = "Huisartsen Assen-West"
$data[1,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsen-assen-west-assen-10018454"
$data[2,0] = "huisartsenpraktijk"
$data[2,1] = "Huisartspraktijk W. van Breugel"
$data[2,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartspraktijk-w-van-breugel-zevenhuizen-3035708"
$data[3,0] = "huisartsenpraktijk"
$data[3,1] = "Huisartspraktijk Olieslagers"
$data[3,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartspraktijk-olieslagers-rockanje-118460"
$data[4,0] = "huisartsenpraktijk"
$data[4,1] = "Huisartsenpraktijk Arts en Zorg Leeuwarden"
$data[4,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-arts-en-zorg-leeuwarden-leeuwarden-3031757"
$data[5,0] = "huisartsenpraktijk"
$data[5,1] = "Huisartsenpraktijk Arts en Zorg Veldweg"
$data[5,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-arts-en-zorg-veldweg-wezep-124911"
$data[6,0] = "huisartsenpraktijk"
$data[6,1] = "Huisartsenpraktijk Medi-Mere Buiten"
$data[6,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-medi-mere-buiten-almere-3031382"
$data[7,0] = "huisartsenpraktijk"
$data[7,1] = "Huisartsenpraktijk Arts en Zorg Gouden Hart"
$data[7,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-arts-en-zorg-gouden-hart-berkel-en-rodenrijs-117268"
$data[8,0] = "huisartsenpraktijk"
$data[8,1] = "Huisartsenpraktijk Arts en Zorg Jan Hendrik"
$data[8,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-arts-en-zorg-jan-hendrik-den-haag-117092"
$data[9,0] = "huisartsenpraktijk"
$data[9,1] = "Gezondheidscentrum Arts en Zorg Winschoten"
$data[9,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-gezondheidscentrum-arts-en-zorg-winschoten-winschoten-3034084"
$data[10,0] = "huisartsenpraktijk"
$data[10,1] = "Huisartsenpraktijk Arts en Zorg Haagse Hout"
$data[10,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-arts-en-zorg-haagse-hout-den-haag-3057129"
$data[11,0] = "huisartsenpraktijk"
$data[11,1] = "Arts en Zorg"
$data[11,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-arts-en-zorg-utrecht-10001581"
$data[12,0] = "huisartsenpraktijk"
$data[12,1] = "Huisartsenpraktijken Medi-Mere"
$data[12,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijken-medi-mere-almere-10026161"
$data[13,0] = "huisartsenpraktijk"
$data[13,1] = "Gezondheidscentrum Arts en Zorg Goudenregenhof"
$data[13,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-gezondheidscentrum-arts-en-zorg-goudenregenhof-den-haag-3006556"
$data[14,0] = "huisartsenpraktijk"
$data[14,1] = "Gezondheidscentrum Arts en Zorg Hoendiep"
$data[14,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-gezondheidscentrum-arts-en-zorg-hoendiep-groningen-251129"
$data[15,0] = "huisartsenpraktijk"
$data[15,1] = "Huisartsenpraktijk Arts en Zorg Vermeertoren"
$data[15,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-arts-en-zorg-vermeertoren-delft-117171"
$data[16,0] = "huisartsenpraktijk"
$data[16,1] = "Fonkelzorg Den Bosch"
$data[16,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-fonkelzorg-den-bosch-s-hertogenbosch-10022918"
$data[17,0] = "huisartsenpraktijk"
$data[17,1] = "Huisartspraktijk Waale"
$data[17,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartspraktijk-waale-breda-205710"
$data[18,0] = "huisartsenpraktijk"
$data[18,1] = "Zorggroep Almere, Huisartsen"
$data[18,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-zorggroep-almere-huisartsen-almere-10021227"
$data[19,0] = "huisartsenpraktijk"
$data[19,1] = "Huisartsenpraktijk A.A. van der Vaart"
$data[19,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-a-a-van-der-vaart-den-haag-3027638"
$data[20,0] = "tandartsenpraktijk"
$data[20,1] = "Tandarts A. Karic-Linic"
$data[20,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandarts-a-karic-linic-den-haag-3042897"
$data[21,0] = "tandartsenpraktijk"
$data[21,1] = "Tandartsjordaan.nl"
$data[21,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartsjordaan-nl-amsterdam-3042877"
$data[22,0] = "tandartsenpraktijk"
$data[22,1] = "Omnident, Kliniek voor Tandheelkunde"
$data[22,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-omnident-kliniek-voor-tandheelkunde-breda-237106"
$data[23,0] = "tandartsenpraktijk"
$data[23,1] = "TandAnders"
$data[23,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandanders-raamsdonksveer-10001160"
$data[24,0] = "tandartsenpraktijk"
$data[24,1] = "Tandartspraktijk Biddinghuizen"
$data[24,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-biddinghuizen-biddinghuizen-3056524"
$data[25,0] = "tandartsenpraktijk"
$data[25,1] = "Tandartspraktijk Waldent"
$data[25,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-waldent-den-haag-10002884"
$data[26,0] = "tandartsenpraktijk"
$data[26,1] = "Tandarts De Ronde Venen"
$data[26,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandarts-de-ronde-venen-mijdrecht-10025363"
$data[27,0] = "tandartsenpraktijk"
$data[27,1] = "Allemans Tandartsen"
$data[27,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-allemans-tandartsen-amerongen-10016011"
$data[28,0] = "tandartsenpraktijk"
$data[28,1] = "Tandheelkundig Centrum Kethel"
$data[28,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandheelkundig-centrum-kethel-schiedam-3055138"
$data[29,0] = "tandartsenpraktijk"
$data[29,1] = "Tandartspraktijk Den Haag Centrum"
$data[29,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-den-haag-centrum-den-haag-3045011"
$data[30,0] = "tandartsenpraktijk"
$data[30,1] = "Tandheelkundig Centrum Wilhelminapier"
$data[30,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandheelkundig-centrum-wilhelminapier-rotterdam-181928"
$data[31,0] = "tandartsenpraktijk"
$data[31,1] = "Tandartspraktijk De Weidevogel"
$data[31,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-de-weidevogel-den-haag-3042694"
$data[32,0] = "tandartsenpraktijk"
$data[32,1] = "SensaDent Tandartsen"
$data[32,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-sensadent-tandartsen-amsterdam-3049885"
$data[33,0] = "tandartsenpraktijk"
$data[33,1] = "TPR | Tandartsenpraktijk Roelofarendsveen"
$data[33,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tpr-tandartsenpraktijk-roelofarendsveen-roelofarendsveen-10017690"
$data[34,0] = "tandartsenpraktijk"
$data[34,1] = "Tandheelkundig Centrum Monnickendam, locatie De Haven"
$data[34,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandheelkundig-centrum-monnickendam-locatie-de-haven-monnickendam-3042997"
$data[35,0] = "tandartsenpraktijk"
$data[35,1] = "Tandartspraktijk Kaptein en Hooykaas"
$data[35,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-kaptein-en-hooykaas-hilversum-179678"
$data[36,0] = "tandartsenpraktijk"
$data[36,1] = "Tandartspraktijk Bas Hengeveld"
$data[36,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-bas-hengeveld-weert-3033605"
$data[37,0] = "tandartsenpraktijk"
$data[37,1] = "Dental Clinics Beesd"
$data[37,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-dental-clinics-beesd-beesd-208613"
$data[38,0] = "tandartsenpraktijk"
$data[38,1] = "De Schans Tandartsen"
$data[38,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-de-schans-tandartsen-leiden-10000705"
$data[39,0] = "tandartsenpraktijk"
$data[39,1] = "Tandartspraktijk Kirsten Knetsch"
$data[39,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-kirsten-knetsch-haarlem-178563"
$data[40,0] = "huisartsenpraktijk"
$data[40,1] = "Huisartsenpraktijk Haverkamp/Bastick"
$data[40,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-haverkamp-bastick-oosterhout-10005075"
$data[41,0] = "huisartsenpraktijk"
$data[41,1] = "Dokters van Hier"
$data[41,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-dokters-van-hier-maastricht-10013414"
$data[42,0] = "huisartsenpraktijk"
$data[42,1] = "Huisartsenpraktijk Kasbergen"
$data[42,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-kasbergen-lunteren-208193"
$data[43,0] = "huisartsenpraktijk"
$data[43,1] = "Huisartsenpraktijk van Elsen"
$data[43,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-van-elsen-oisterwijk-10019665"
$data[44,0] = "huisartsenpraktijk"
$data[44,1] = "Apotheekhoudende huisartsenpraktijk J.E. de Groot"
$data[44,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-apotheekhoudende-huisartsenpraktijk-j-e-de-groot-zuid-beijerland-3060002"
$data[45,0] = "huisartsenpraktijk"
$data[45,1] = "Huisartsenpraktijk Boddeus en Steenbergen"
$data[45,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-boddeus-en-steenbergen-groningen-126018"
$data[46,0] = "huisartsenpraktijk"
$data[46,1] = "Huisartspraktijk Van Eijk"
$data[46,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartspraktijk-van-eijk-delft-207863"
$data[47,0] = "huisartsenpraktijk"
$data[47,1] = "Huisartsen Kievit"
$data[47,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsen-kievit-bussum-205922"
$data[48,0] = "huisartsenpraktijk"
$data[48,1] = "Huisartsenpraktijk Binck-Zorg"
$data[48,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-binck-zorg-den-haag-10020425"
$data[49,0] = "huisartsenpraktijk"
$data[49,1] = "Huisartsenpraktijk Boas & Valkenburg"
$data[49,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-boas-valkenburg-vianen-10012743"
$data[50,0] = "huisartsenpraktijk"
$data[50,1] = "Huisartsenpraktijk Kemper en Oldenburg"
$data[50,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-kemper-en-oldenburg-alkmaar-206554"
$data[51,0] = "huisartsenpraktijk"
$data[51,1] = "Huisartsenpraktijk Ridha"
$data[51,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-ridha-volendam-3048383"
$data[52,0] = "huisartsenpraktijk"
$data[52,1] = "Huisartsenpraktijk Trompert"
$data[52,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-trompert-zuidhorn-207393"
$data[53,0] = "huisartsenpraktijk"
$data[53,1] = "Huisartsenpraktijk Uitvindersbuurt"
$data[53,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-uitvindersbuurt-ede-123316"
$data[54,0] = "huisartsenpraktijk"
$data[54,1] = "ClydesCure, Huisartsenpraktijk"
$data[54,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-clydescure-huisartsenpraktijk-lelystad-10082191"
$data[55,0] = "huisartsenpraktijk"
$data[55,1] = "Huisartsenpraktijk Blessing - Jeuring - Moonen"
$data[55,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-blessing-jeuring-moonen-udenhout-120894"
$data[56,0] = "huisartsenpraktijk"
$data[56,1] = "Huisartsenpraktijk Bongers"
$data[56,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-bongers-zeewolde-207820"
$data[57,0] = "huisartsenpraktijk"
$data[57,1] = "Huisartsenpraktijk De Diependaal"
$data[57,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-de-diependaal-stein-237109"
$data[58,0] = "huisartsenpraktijk"
$data[58,1] = "Huisartsenpraktijk De Dolfijn"
$data[58,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-de-dolfijn-middelburg-120118"
$data[59,0] = "huisartsenpraktijk"
$data[59,1] = "Huisartsenpraktijk De Hooge Boom"
$data[59,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-de-hooge-boom-hoogwoud-3005987"
$data[60,0] = "tandartsenpraktijk"
$data[60,1] = "Ferguson Hannewijk Tandartsen"
$data[60,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-ferguson-hannewijk-tandartsen-rijnsburg-180935"
$data[61,0] = "tandartsenpraktijk"
$data[61,1] = "Tandartsenpraktijk Zoeterwoudsesingel"
$data[61,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartsenpraktijk-zoeterwoudsesingel-leiden-180956"
$data[62,0] = "tandartsenpraktijk"
$data[62,1] = "Tandheelkunde & Implantologie Amsterdam"
$data[62,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandheelkunde-implantologie-amsterdam-amsterdam-10026757"
$data[63,0] = "tandartsenpraktijk"
$data[63,1] = "Tandheelkunde & Implantologie Amsterdam, locatie Stadionweg"
$data[63,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandheelkunde-implantologie-amsterdam-locatie-stadionweg-amsterdam-178713"
$data[64,0] = "tandartsenpraktijk"
$data[64,1] = "Tandartspraktijk Dronten"
$data[64,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-dronten-dronten-10017070"
$data[65,0] = "tandartsenpraktijk"
$data[65,1] = "Tandartspraktijk Inge Schrauwen"
$data[65,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-inge-schrauwen-tilburg-3033485"
$data[66,0] = "tandartsenpraktijk"
$data[66,1] = "Dental Clinics Zandvoort"
$data[66,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-dental-clinics-zandvoort-zandvoort-180752"
$data[67,0] = "tandartsenpraktijk"
$data[67,1] = "Mondzorgcentrum Takenhofplein"
$data[67,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-mondzorgcentrum-takenhofplein-nijmegen-3059439"
$data[68,0] = "tandartsenpraktijk"
$data[68,1] = "Praktijk voor tandheelkunde en mondhygiëne Van Uijtert"
$data[68,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-praktijk-voor-tandheelkunde-en-mondhygiene-van-uijtert-lage-zwaluwe-3035680"
$data[69,0] = "tandartsenpraktijk"
$data[69,1] = "Dental Clinics Purmerend de Gors"
$data[69,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-dental-clinics-purmerend-de-gors-purmerend-179617"
$data[70,0] = "tandartsenpraktijk"
$data[70,1] = "KTA Kliniek voor Tandheelkunde Amersfoort"
$data[70,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-kta-kliniek-voor-tandheelkunde-amersfoort-amersfoort-3024325"
$data[71,0] = "tandartsenpraktijk"
$data[71,1] = "Tandartspraktijk De Jol Lelystad"
$data[71,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-de-jol-lelystad-lelystad-3056450"
$data[72,0] = "tandartsenpraktijk"
$data[72,1] = "DentXperts"
$data[72,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-dentxperts-amsterdam-10005995"
$data[73,0] = "tandartsenpraktijk"
$data[73,1] = "Tandheelkundig Centrum Nederland Emmastraat"
$data[73,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandheelkundig-centrum-nederland-emmastraat-alkmaar-180307"
$data[74,0] = "tandartsenpraktijk"
$data[74,1] = "Mondzorg Aveling"
$data[74,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-mondzorg-aveling-hoogvliet-10003794"
$data[75,0] = "tandartsenpraktijk"
$data[75,1] = "Tandartspraktijk Claessens"
$data[75,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartspraktijk-claessens-echt-185745"
$data[76,0] = "tandartsenpraktijk"
$data[76,1] = "Tandheelkundig Centrum Koraalzwam"
$data[76,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandheelkundig-centrum-koraalzwam-alphen-aan-den-rijn-3041017"
$data[77,0] = "tandartsenpraktijk"
$data[77,1] = "Mondzorg Maarssen"
$data[77,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-mondzorg-maarssen-maarssen-179758"
$data[78,0] = "tandartsenpraktijk"
$data[78,1] = "Tandartsenpraktijk Steyl"
$data[78,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-tandartsenpraktijk-steyl-steyl-185551"
$data[79,0] = "tandartsenpraktijk"
$data[79,1] = "Top Dental"
$data[79,2] = "https://www.zorgkaartnederland.nl/zorginstelling/tandartsenpraktijk-top-dental-volendam-3049583"
$data[80,0] = "huisartsenpraktijk"
$data[80,1] = "Huisartsenpraktijk Havekes"
$data[80,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-havekes-harderwijk-10019799"
$data[81,0] = "huisartsenpraktijk"
$data[81,1] = "Huisartsenpraktijk Het Zorgkwartier"
$data[81,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-het-zorgkwartier-ommen-124533"
$data[82,0] = "huisartsenpraktijk"
$data[82,1] = "Huisartsenpraktijk J.P. van Tussenbroek"
$data[82,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-j-p-van-tussenbroek-delft-117233"
$data[83,0] = "huisartsenpraktijk"
$data[83,1] = "Huisartsenpraktijk Majdandzic"
$data[83,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-majdandzic-breda-205739"
$data[84,0] = "huisartsenpraktijk"
$data[84,1] = "Huisartsenpraktijk Ritter en Kuipers"
$data[84,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-ritter-en-kuipers-rotterdam-208445"
$data[85,0] = "huisartsenpraktijk"
$data[85,1] = "Huisartsenpraktijk Schoterpoort, Praktijk Steketee"
$data[85,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-schoterpoort-praktijk-steketee-haarlem-3046049"
$data[86,0] = "huisartsenpraktijk"
$data[86,1] = "Huisartsenpraktijk Tabak"
$data[86,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-tabak-amsterdam-3046193"
$data[87,0] = "huisartsenpraktijk"
$data[87,1] = "Huisartsenpraktijk Ulestraten"
$data[87,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-ulestraten-ulestraten-3031904"
$data[88,0] = "huisartsenpraktijk"
$data[88,1] = "Huisartsenpraktijk van Beijsterveldt, locatie Dorst"
$data[88,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-van-beijsterveldt-locatie-dorst-dorst-120647"
$data[89,0] = "huisartsenpraktijk"
$data[89,1] = "Huisartsenpraktijk Vlaslant, Praktijk Verhoeckx"
$data[89,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-vlaslant-praktijk-verhoeckx-valkenswaard-10010787"
$data[90,0] = "huisartsenpraktijk"
$data[90,1] = "Huisartsenpraktijk Willems"
$data[90,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-willems-best-3056625"
$data[91,0] = "huisartsenpraktijk"
$data[91,1] = "Huisartsenpraktijk ZorghoekWestland"
$data[91,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartsenpraktijk-zorghoekwestland-honselersdijk-10001572"
$data[92,0] = "huisartsenpraktijk"
$data[92,1] = "Huisartspraktijk L. van Eijk"
$data[92,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartspraktijk-l-van-eijk-den-dolder-3037015"
$data[93,0] = "huisartsenpraktijk"
$data[93,1] = "Huisartspraktijk Willemsen"
$data[93,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-huisartspraktijk-willemsen-capelle-aan-den-ijssel-117569"
$data[94,0] = "huisartsenpraktijk"
$data[94,1] = "Apotheekhoudende Huisartsenpraktijk De Krim"
$data[94,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-apotheekhoudende-huisartsenpraktijk-de-krim-de-krim-124570"
$data[95,0] = "huisartsenpraktijk"
$data[95,1] = "Apotheekhoudende Huisartsenpraktijk Feij en Van der Wal"
$data[95,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-apotheekhoudende-huisartsenpraktijk-feij-en-van-der-wal-slochteren-3046266"
$data[96,0] = "huisartsenpraktijk"
$data[96,1] = "Apotheekhoudende huisartspraktijk J. de Kroon"
$data[96,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-apotheekhoudende-huisartspraktijk-j-de-kroon-onstwedde-207076"
$data[97,0] = "huisartsenpraktijk"
$data[97,1] = "De Bergense Huisartsen- Praktijk Visser"
$data[97,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-de-bergense-huisartsen-praktijk-visser-bergen-3043935"
$data[98,0] = "huisartsenpraktijk"
$data[98,1] = "De Colvenier, Huisartspraktijk D.H.A. Pons"
$data[98,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-de-colvenier-huisartspraktijk-d-h-a-pons-gorinchem-120010"
$data[99,0] = "huisartsenpraktijk"
$data[99,1] = "GCM Gezondheidscentrum Boomstede, Huisartsen"
$data[99,2] = "https://www.zorgkaartnederland.nl/zorginstelling/huisartsenpraktijk-gcm-gezondheidscentrum-boomstede-huisartsen-maarssen-102611"

# Rows 1-21 (header + first 20 tandartsenpraktijk entries) stay untouched;
# rows 22-121 are written in one shot to reflect the new, longer dataset.
$ws.Range("A22:C121").Value2 = $data
